$wb = $excel.ActiveWorkbook

# The shared string "Ready for handoff" (used by the Overview sheet's
# zh-cn/de-de status cells for the cae33588... row, and by the "Status"
# column of both the zh-cn and de-de detail sheets for that same row) is
# replaced throughout with "Handback transform failed".
$newStatus = "Handback transform failed"

# --- Overview sheet: row 3 (cae33588... file) status for both locales.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status column (C) for row 3, plus populate the new
# --- "Error Detail" column (P) for row 3 and widen the column to fit it.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("P3").Value = "Handback file name: aeboa3l5.1d3 is different with handoff file name: cae33588-4425-4cc3-9990-5bbeeb7febff.c89148643c39e0f32dc38af26446ba98888bdf01.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: same Status + Error Detail population + column width.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("P3").Value = "Handback file name: aeboa3l5.1d3 is different with handoff file name: cae33588-4425-4cc3-9990-5bbeeb7febff.c89148643c39e0f32dc38af26446ba98888bdf01.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
